$d = $word.ActiveDocument

# --- Paragraph 1: add a paragraph border, widen left indent, and collapse
#     the "**ID__..." run + trailing lone-space run into a single run with
#     the updated bookmark-style id text. ---
$p1 = $d.Paragraphs(1)

# Add pBdr (top/left/bottom/right, 5pt space each)
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# Remove the trailing run that contains only a single space, then update
# the remaining run's text in place (so the two runs collapse into one).
$p1Range = $p1.Range
$spaceRange = $d.Range($p1Range.End - 2, $p1Range.End - 1)
$spaceRange.Delete()

$idRange = $d.Range($p1Range.Start, $p1Range.End - 1)
$idRange.Text = "**ID__AFFARS_AF_PGI_5349_101__ID**"
